$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (date), Volumen, Precio máximo, Precio promedio ponderado
# and Precio $/Kg values between row 2 and row 3.

$ws.Range("D2").Value = 44672
$ws.Range("M2").Value = 8
$ws.Range("O2").Value = 180000
$ws.Range("P2").Value = 180000
$ws.Range("S2").Value = 180000

$ws.Range("D3").Value = 44253
$ws.Range("M3").Value = 12
$ws.Range("O3").Value = 200000
$ws.Range("P3").Value = 190000
$ws.Range("S3").Value = 190000
